$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.230.77'
$ws.Range('E2').Value = '  +1.03%  '

$ws.Range('D3').Value = '3.120.39'
$ws.Range('E3').Value = '  +1.88%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '398.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.30%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.539'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.12%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.08%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.598'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.82%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.92'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.89%  '

$ws.Range('E11').Value = '  +0.81%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0860'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.62%  '

$ws.Range('D13').Value = '3.605.47'
$ws.Range('E13').Value = '  +1.52%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.82'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.36%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.82'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.52%  '

$ws.Range('E16').Value = '  +5.93%  '

$ws.Range('D17').Value = '3.124.64'
$ws.Range('E17').Value = '  +1.83%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.23%  '

$ws.Range('D19').Value = '52.056.62'
$ws.Range('E19').Value = '  +0.54%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.42%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.98%  '

$ws.Range('D22').Value = '0.0₃0972'
$ws.Range('E22').Value = '  +0.66%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.74%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.40%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.28%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.75%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.33'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.67%  '

$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.168'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.88%  '

$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.02%  '

$ws.Range('E31').Value = '  -0.61%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.90'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.64%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0490'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.22%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '36.36'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.91%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.08'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.62%  '

$ws.Range('E36').Value = '  -0.88%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.43'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.62%  '

$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.294'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.21%  '

$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.06'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.85%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.66'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.44%  '

$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.75%  '

$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '130.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.31%  '

$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.54%  '

$ws.Range('E45').Value = '  -0.13%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.15'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.23%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.06'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.48%  '

$ws.Range('D49').Value = '2.089.97'
$ws.Range('E49').Value = '  +1.87%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0521'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +32.84%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.930'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.63%  '
